$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BF2:BF31 hold the game "Date" column. It was previously populated with
# the literal text "5-6-2007-08" (an artifact of the folder/file naming
# convention used for the season), but the NBA box-score this data was
# scraped from actually reports the game date as May 6, 2008 local time -
# the stats were pulled a day off from how the NBA site displayed them.
# Correct all 30 data rows to the text "2008-05-06".
#
# NumberFormat is forced to Text ("@") before the assignment so Excel's
# automatic type inference doesn't silently turn the literal string into
# a date serial number; the format is then reset back to the default
# "Normal" cell style so no visible formatting change is introduced.
$rng = $ws.Range("BF2:BF31")
$rng.NumberFormat = "@"
$rng.Value = "2008-05-06"
$rng.Style = "Normal"
